# Applies the "Short tutorial documents start" edit:
#  1. NCPResultsV1        -> NCPResults
#  2. NCPSettingsV1.xml   -> NCPSettings.xml   (two occurrences)
#  3. "Run NCP V1"        -> "Run NCP"         (drop the " V1" run)
#  4. Remove the "Experiment with different numbers of synergies" and
#     "Experiment with activation minimization" sections (everything
#     after the "...Inverse Dynamics joint moments" paragraph), leaving
#     that paragraph terminated with a period.

$d = $word.ActiveDocument

# --- 4. Truncate the tutorial -------------------------------------------
# Paragraph 44 ("Plot 3 ... Inverse Dynamics joint moments") is the last
# paragraph we keep; everything from paragraph 45 through the final
# paragraph gets removed in one shot, then we append the closing period.
$pCut = $d.Paragraphs.Item(45)
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$cutRange = $d.Range($pCut.Range.Start, $pLast.Range.End)
$cutRange.Delete()

$pFinal = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($pFinal.Range.End - 1, $pFinal.Range.End - 1)
$insertPoint.Text = "."

# --- 1. NCPResultsV1 -> NCPResults ---------------------------------------
$d.Content.Find.Execute("NCPResultsV1", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "NCPResults", 2)

# --- 2. NCPSettingsV1.xml -> NCPSettings.xml -----------------------------
$d.Content.Find.Execute("NCPSettingsV1.xml", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "NCPSettings.xml", 2)

# --- 3. "Run NCP V1" -> "Run NCP" ----------------------------------------
$d.Content.Find.Execute("Run NCP V1", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Run NCP", 2)
